$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 349.5
$ws.Range("I4").Value = 199.5
$ws.Range("J4").Value = 499.5
$ws.Range("K4").Value = 199.5
$ws.Range("L4").Value = 499.5
$ws.Range("M4").Value = -85.5
$ws.Range("N4").Value = -727.5
$ws.Range("H92").Value = 37037704
$ws.Range("I92").Value = 37037704
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 37037704
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -37036456
$ws.Range("N92").ClearContents()
$ws.Range("H138").Value = 3984761.8
$ws.Range("I138").Value = 1610878.9
$ws.Range("J138").Value = 4764752
$ws.Range("K138").Value = 4832636.699999999
$ws.Range("L138").Value = 14294256
$ws.Range("M138").Value = -4827496.699999999
$ws.Range("N138").Value = -14304536

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2626.7727
$ws.Range("I61").Value = 1721.6666
$ws.Range("J61").Value = 4566.2856
$ws.Range("K61").Value = 1721.6666
$ws.Range("L61").Value = 4566.2856
$ws.Range("M61").Value = -1509.6666
$ws.Range("N61").Value = -4990.2856
$ws.Range("H88").Value = 7999.6665
$ws.Range("I88").Value = 1999.5
$ws.Range("K88").Value = 1999.5
$ws.Range("M88").Value = -1593.5
$ws.Range("H91").Value = 7999.6665
$ws.Range("I91").Value = 1999.5
$ws.Range("K91").Value = 1999.5
$ws.Range("M91").Value = -595.5
$ws.Range("H122").Value = 1603.1892
$ws.Range("I122").Value = 1312.24
$ws.Range("K122").Value = 3936.72
$ws.Range("M122").Value = -1486.72
$ws.Range("H136").Value = 2626.7727
$ws.Range("I136").Value = 1721.6666
$ws.Range("J136").Value = 4566.2856
$ws.Range("K136").Value = 5164.9998
$ws.Range("L136").Value = 13698.8568
$ws.Range("M136").Value = -2614.9998
$ws.Range("N136").Value = -18798.8568

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5295.407
$ws.Range("I86").Value = 1870.7693
$ws.Range("J86").Value = 8475.429
$ws.Range("K86").Value = 1870.7693
$ws.Range("L86").Value = 8475.429
$ws.Range("M86").Value = -747.7692999999999
$ws.Range("N86").Value = -10721.429
$ws.Range("H89").Value = 5295.407
$ws.Range("I89").Value = 1870.7693
$ws.Range("J89").Value = 8475.429
$ws.Range("K89").Value = 9353.8465
$ws.Range("L89").Value = 42377.145
$ws.Range("M89").Value = -3737.8465
$ws.Range("N89").Value = -53609.145
$ws.Range("H105").Value = 3366.6333
$ws.Range("I105").Value = 3136.318
$ws.Range("K105").Value = 3136.318
$ws.Range("M105").Value = -1389.318

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 18324.666
$ws.Range("I33").Value = 18324.666
$ws.Range("K33").Value = 18324.666
$ws.Range("M33").Value = -17945.666
$ws.Range("H41").Value = 11499.5
$ws.Range("I41").Value = 11499.5
$ws.Range("K41").Value = 11499.5
$ws.Range("M41").Value = -11071.5
$ws.Range("H58").Value = 2350.84
$ws.Range("I58").Value = 1452.6154
$ws.Range("J58").Value = 3323.9167
$ws.Range("K58").Value = 1452.6154
$ws.Range("L58").Value = 3323.9167
$ws.Range("M58").Value = -1249.6154
$ws.Range("N58").Value = -3729.9167
$ws.Range("H62").Value = 32643.428
$ws.Range("I62").Value = 52626.25
$ws.Range("J62").Value = 5999.6665
$ws.Range("K62").Value = 52626.25
$ws.Range("L62").Value = 5999.6665
$ws.Range("M62").Value = -52002.25
$ws.Range("N62").Value = -7247.6665
$ws.Range("H65").Value = 32643.428
$ws.Range("I65").Value = 52626.25
$ws.Range("J65").Value = 5999.6665
$ws.Range("K65").Value = 263131.25
$ws.Range("L65").Value = 29998.3325
$ws.Range("M65").Value = -260011.25
$ws.Range("N65").Value = -36238.3325
$ws.Range("H98").Value = 36999
$ws.Range("J98").Value = 36999
$ws.Range("L98").Value = 36999
$ws.Range("N98").Value = -41491
$ws.Range("H99").Value = 8929629
$ws.Range("I99").Value = 8929629
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 8929629
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -8928131
$ws.Range("N99").ClearContents()
$ws.Range("H122").Value = 1966.9286
$ws.Range("I122").Value = 980.2222
$ws.Range("J122").Value = 3743
$ws.Range("K122").Value = 2940.6666
$ws.Range("L122").Value = 11229
$ws.Range("M122").Value = -490.6666
$ws.Range("N122").Value = -16129
$ws.Range("H126").Value = 8929629
$ws.Range("I126").Value = 8929629
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 26788887
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -26786417
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 2350.84
$ws.Range("I136").Value = 1452.6154
$ws.Range("J136").Value = 3323.9167
$ws.Range("K136").Value = 4357.8462
$ws.Range("L136").Value = 9971.750100000001
$ws.Range("M136").Value = -1807.8462
$ws.Range("N136").Value = -15071.7501

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 960
$ws.Range("I51").Value = 950
$ws.Range("J51").Value = 966.6667
$ws.Range("K51").Value = 2850
$ws.Range("L51").Value = 2900.0001
$ws.Range("M51").Value = -2390
$ws.Range("N51").Value = -3820.0001
$ws.Range("H117").Value = 973.75
$ws.Range("J117").Value = 877.3333
$ws.Range("L117").Value = 2631.9999
$ws.Range("N117").Value = -9515.999899999999
$ws.Range("H119").Value = 1103.1666
$ws.Range("I119").Value = 254.75
$ws.Range("J119").Value = 2800
$ws.Range("K119").Value = 764.25
$ws.Range("L119").Value = 8400
$ws.Range("M119").Value = 4073.75
$ws.Range("N119").Value = -18076
$ws.Range("H129").Value = 1283.875
$ws.Range("J129").Value = 1722.2
$ws.Range("L129").Value = 5166.6
$ws.Range("N129").Value = -15166.6
$ws.Range("H131").Value = 1322.2162
$ws.Range("I131").Value = 337.53845
$ws.Range("J131").Value = 1532.0656
$ws.Range("K131").Value = 1012.61535
$ws.Range("L131").Value = 4596.1968
$ws.Range("M131").Value = 4027.38465
$ws.Range("N131").Value = -14676.1968
$ws.Range("H132").Value = 1245.6316
$ws.Range("I132").Value = 770.8
$ws.Range("J132").Value = 1415.2142
$ws.Range("K132").Value = 6937.2
$ws.Range("L132").Value = 12736.9278
$ws.Range("M132").Value = -4407.2
$ws.Range("N132").Value = -17796.9278

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5364.3555
$ws.Range("J70").Value = 5355
$ws.Range("L70").Value = 5355
$ws.Range("N70").Value = -5895
$ws.Range("H73").Value = 5364.3555
$ws.Range("J73").Value = 5355
$ws.Range("L73").Value = 5355
$ws.Range("N73").Value = -7227

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2974.1304
$ws.Range("I7").Value = 1757.1428
$ws.Range("K7").Value = 1757.1428
$ws.Range("M7").Value = -1645.1428
$ws.Range("H22").Value = 1142
$ws.Range("I22").Value = 400.66666
$ws.Range("J22").Value = 1883.3334
$ws.Range("K22").Value = 400.66666
$ws.Range("L22").Value = 1883.3334
$ws.Range("M22").Value = -105.66666
$ws.Range("N22").Value = -2473.3334
$ws.Range("H27").Value = 1142
$ws.Range("I27").Value = 400.66666
$ws.Range("J27").Value = 1883.3334
$ws.Range("K27").Value = 400.66666
$ws.Range("L27").Value = 1883.3334
$ws.Range("M27").Value = -293.66666
$ws.Range("N27").Value = -2097.3334
$ws.Range("H40").Value = 2333.853
$ws.Range("I40").Value = 1227.9286
$ws.Range("J40").Value = 3108
$ws.Range("K40").Value = 1227.9286
$ws.Range("L40").Value = 3108
$ws.Range("M40").Value = -1091.9286
$ws.Range("N40").Value = -3380
$ws.Range("H93").Value = 1364.6
$ws.Range("I93").Value = 1057.4286
$ws.Range("J93").Value = 2081.3333
$ws.Range("K93").Value = 1057.4286
$ws.Range("L93").Value = 2081.3333
$ws.Range("M93").Value = 190.5714
$ws.Range("N93").Value = -4577.3333
$ws.Range("H126").Value = 2974.1304
$ws.Range("I126").Value = 1757.1428
$ws.Range("K126").Value = 5271.428400000001
$ws.Range("M126").Value = -2801.428400000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3473104.5
$ws.Range("I107").Value = 5051426.5
$ws.Range("J107").Value = 796.6
$ws.Range("K107").Value = 15154279.5
$ws.Range("L107").Value = 2389.8
$ws.Range("M107").Value = -15152359.5
$ws.Range("N107").Value = -6229.8
$ws.Range("H122").Value = 28529.236
$ws.Range("I122").Value = 41684.36
$ws.Range("J122").Value = 3230.923
$ws.Range("K122").Value = 125053.08
$ws.Range("L122").Value = 9692.769
$ws.Range("M122").Value = -122603.08
$ws.Range("N122").Value = -14592.769
$ws.Range("H126").Value = 47495.363
$ws.Range("I126").Value = 60924.59
$ws.Range("J126").Value = 1836
$ws.Range("K126").Value = 182773.77
$ws.Range("L126").Value = 5508
$ws.Range("M126").Value = -180303.77
$ws.Range("N126").Value = -10448
$ws.Range("H136").Value = 7269476
$ws.Range("I136").Value = 9287938
$ws.Range("J136").Value = 3013.4
$ws.Range("K136").Value = 27863814
$ws.Range("L136").Value = 9040.2
$ws.Range("M136").Value = -27861264
$ws.Range("N136").Value = -14140.2
